# Refresh Gilgamesh Profits market-data snapshot (scheduled runner sync).
# For each affected Leve row, currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are updated to the latest Universalis pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 589.2727
$ws.Range("I2").Value = 360
$ws.Range("J2").Value = 864.4
$ws.Range("K2").Value = 360
$ws.Range("L2").Value = 864.4
$ws.Range("M2").Value = -247
$ws.Range("N2").Value = -1090.4

# Row 18
$ws.Range("H18").Value = 47625276
$ws.Range("I18").Value = 6374.8125
$ws.Range("J18").Value = 200005760
$ws.Range("K18").Value = 6374.8125
$ws.Range("L18").Value = 200005760
$ws.Range("M18").Value = -6090.8125
$ws.Range("N18").Value = -200006328

# Row 40
$ws.Range("H40").Value = 5220.5557
$ws.Range("I40").Value = 4848.3335
$ws.Range("K40").Value = 4848.3335
$ws.Range("M40").Value = -4673.3335

# Row 55
$ws.Range("H55").Value = 508
$ws.Range("J55").Value = 397.5
$ws.Range("L55").Value = 397.5
$ws.Range("N55").Value = -825.5

# Row 64
$ws.Range("H64").Value = 250006990
$ws.Range("I64").Value = 9333
$ws.Range("K64").Value = 9333
$ws.Range("M64").Value = -9085

# Row 67
$ws.Range("H67").Value = 250006990
$ws.Range("I67").Value = 9333
$ws.Range("K67").Value = 9333
$ws.Range("M67").Value = -8475

# Row 103
$ws.Range("H103").Value = 2644.2727
$ws.Range("I103").Value = 3181.8333
$ws.Range("K103").Value = 9545.499899999999
$ws.Range("M103").Value = -8959.499899999999

# Row 113
$ws.Range("H113").Value = 1196
$ws.Range("J113").Value = 1196
$ws.Range("L113").Value = 1196
$ws.Range("N113").Value = -7704

# Row 137
$ws.Range("H137").Value = 3856.9773
$ws.Range("I137").Value = 1529.6
$ws.Range("J137").Value = 8844.214
$ws.Range("K137").Value = 4588.799999999999
$ws.Range("L137").Value = 26532.642
$ws.Range("M137").Value = -2038.799999999999
$ws.Range("N137").Value = -31632.642

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 175992.1
$ws.Range("I74").Value = 266273.38
$ws.Range("J74").Value = 3636.9092
$ws.Range("K74").Value = 266273.38
$ws.Range("L74").Value = 3636.9092
$ws.Range("M74").Value = -265399.38
$ws.Range("N74").Value = -5384.9092

# Row 77
$ws.Range("H77").Value = 175992.1
$ws.Range("I77").Value = 266273.38
$ws.Range("J77").Value = 3636.9092
$ws.Range("K77").Value = 1331366.9
$ws.Range("L77").Value = 18184.546
$ws.Range("M77").Value = -1326998.9
$ws.Range("N77").Value = -26920.546

# Row 122
$ws.Range("H122").Value = 1246.1818
$ws.Range("I122").Value = 1088.625
$ws.Range("K122").Value = 3265.875
$ws.Range("M122").Value = -815.875

$ws = $wb.Worksheets.Item("BSM")
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 134
$ws.Range("H134").Value = 2762.725
$ws.Range("I134").Value = 2509.8276
$ws.Range("K134").Value = 7529.4828
$ws.Range("M134").Value = -4994.4828

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3905.6155
$ws.Range("I31").Value = 3114.8076
$ws.Range("K31").Value = 3114.8076
$ws.Range("M31").Value = -2819.8076

# Row 34
$ws.Range("H34").Value = 3905.6155
$ws.Range("I34").Value = 3114.8076
$ws.Range("K34").Value = 3114.8076
$ws.Range("M34").Value = -2912.8076

# Row 58
$ws.Range("H58").Value = 4433
$ws.Range("I58").Value = 4320.5713
$ws.Range("J58").Value = 4493.5386
$ws.Range("K58").Value = 4320.5713
$ws.Range("L58").Value = 4493.5386
$ws.Range("M58").Value = -4117.5713
$ws.Range("N58").Value = -4899.5386

# Row 115
$ws.Range("H115").Value = 49913.5
$ws.Range("J115").Value = 49913.5
$ws.Range("L115").Value = 49913.5
$ws.Range("N115").Value = -52263.5

# Row 122
$ws.Range("H122").Value = 3656.2144
$ws.Range("I122").Value = 3070.5557
$ws.Range("J122").Value = 4710.4
$ws.Range("K122").Value = 9211.667099999999
$ws.Range("L122").Value = 14131.2
$ws.Range("M122").Value = -6761.667099999999
$ws.Range("N122").Value = -19031.2

# Row 132
$ws.Range("H132").Value = 2359.9153
$ws.Range("I132").Value = 1829.8723
$ws.Range("K132").Value = 5489.6169
$ws.Range("M132").Value = -2959.6169

# Row 136
$ws.Range("H136").Value = 4433
$ws.Range("I136").Value = 4320.5713
$ws.Range("J136").Value = 4493.5386
$ws.Range("K136").Value = 12961.7139
$ws.Range("L136").Value = 13480.6158
$ws.Range("M136").Value = -10411.7139
$ws.Range("N136").Value = -18580.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 137
$ws.Range("H137").Value = 4223.095
$ws.Range("J137").Value = 4109.5835
$ws.Range("L137").Value = 12328.7505
$ws.Range("N137").Value = -22528.7505

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 142859660
$ws.Range("I80").Value = 250001680
$ws.Range("J80").Value = 3649.6667
$ws.Range("K80").Value = 250001680
$ws.Range("L80").Value = 3649.6667
$ws.Range("M80").Value = -250000682
$ws.Range("N80").Value = -5645.6667

# Row 83
$ws.Range("H83").Value = 142859660
$ws.Range("I83").Value = 250001680
$ws.Range("J83").Value = 3649.6667
$ws.Range("K83").Value = 1250008400
$ws.Range("L83").Value = 18248.3335
$ws.Range("M83").Value = -1250003408
$ws.Range("N83").Value = -28232.3335

# Row 122
$ws.Range("H122").Value = 2728.818
$ws.Range("I122").Value = 2780
$ws.Range("J122").Value = 2498.5
$ws.Range("K122").Value = 8340
$ws.Range("L122").Value = 7495.5
$ws.Range("M122").Value = -5890
$ws.Range("N122").Value = -12395.5

# Row 126
$ws.Range("H126").Value = 4255.5264
$ws.Range("J126").Value = 7388
$ws.Range("L126").Value = 22164
$ws.Range("N126").Value = -27104

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6083.0835
$ws.Range("I40").Value = 5913
$ws.Range("K40").Value = 5913
$ws.Range("M40").Value = -5777

# Row 132
$ws.Range("H132").Value = 5222.8076
$ws.Range("I132").Value = 2988.3684
$ws.Range("J132").Value = 11287.714
$ws.Range("K132").Value = 8965.1052
$ws.Range("L132").Value = 33863.142
$ws.Range("M132").Value = -6435.1052
$ws.Range("N132").Value = -38923.142

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2671.4807
$ws.Range("I132").Value = 2825.1892
$ws.Range("J132").Value = 2292.3333
$ws.Range("K132").Value = 8475.567599999998
$ws.Range("L132").Value = 6876.999899999999
$ws.Range("M132").Value = -5945.567599999998
$ws.Range("N132").Value = -11936.9999

